$d = $word.ActiveDocument

function Set-ParagraphXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $pRange = $p.Range
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body>' + $innerXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData>' +
        '</pkg:part>' +
        '</pkg:package>'
    $pRange.InsertXML($xml)
}

# 1. "8 October 2023" -> "28" + " October 2023" (two runs, date corrected to 28 October 2023)
$xml1 = '<w:p w14:paraId="5933949F" w14:textId="4227AD85" w:rsidR="006F26B6" w:rsidRDefault="006F26B6">' + `
    '<w:r><w:t>28</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> October 2023</w:t></w:r>' + `
    '</w:p>'
Set-ParagraphXml 2 $xml1

# 2. Merge "with " + "Flask-" runs into a single "with Flask-" run
$xml2 = '<w:p w14:paraId="780E3B03" w14:textId="6735E58B" w:rsidR="009B5032" w:rsidRPr="00F52218" w:rsidRDefault="009B5032">' + `
    '<w:r><w:t xml:space="preserve">1) </w:t></w:r>' + `
    '<w:r w:rsidRPr="00F52218"><w:t xml:space="preserve">Set up a local database and read data in </w:t></w:r>' + `
    '<w:r w:rsidR="00BF6A96" w:rsidRPr="00F52218"><w:t xml:space="preserve">it </w:t></w:r>' + `
    '<w:r w:rsidRPr="00F52218"><w:t>with Flask-</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r w:rsidRPr="00F52218"><w:t>SQLAlchemy</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '</w:p>'
Set-ParagraphXml 25 $xml2

# 3. Merge "2) " + "Clean your data and perform a " runs into a single run
$xml3 = '<w:p w14:paraId="31734EB3" w14:textId="6E1ED372" w:rsidR="00BF6A96" w:rsidRPr="00F52218" w:rsidRDefault="00BF6A96">' + `
    '<w:r w:rsidRPr="00F52218"><w:t xml:space="preserve">2) Clean your data and perform a </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r w:rsidRPr="00F52218"><w:t>pandas</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r w:rsidRPr="00F52218"><w:t xml:space="preserve"> merge with your two data sets, then calculate some new values based on the new data set.</w:t></w:r>' + `
    '</w:p>'
Set-ParagraphXml 26 $xml3

# 4. Merge " " + "Make at least 1 Pandas..." runs into a single run
$xml4 = '<w:p w14:paraId="0CD25D1F" w14:textId="6D72A2B8" w:rsidR="00BF6A96" w:rsidRDefault="00BF6A96">' + `
    '<w:r w:rsidRPr="00F52218"><w:t>3)</w:t></w:r>' + `
    '<w:r w:rsidR="00F52218" w:rsidRPr="00F52218"><w:t xml:space="preserve"> Make at least 1 Pandas pivot table and 1 matplotlib/seaborn plot. Pivot tables are a way to summarize your data and present it easily in a way that isn&#8217;t just a graph. They can be useful when combined with graphs.</w:t></w:r>' + `
    '</w:p>'
Set-ParagraphXml 27 $xml4
